# #CRM-31 Remove ID, Bank Details, bracket flag from Download SF list
#
# The SF List download template had a few columns that should no longer be
# exported: the vendor "ID" column, the four bank-detail columns (Bank Name,
# Bank Account, IFSC Code, Beneficiary Name) and the "Brackets Flag" column.
# Deleting them (rather than just blanking them) shifts every column to its
# right one step to the left, which is what the target workbook shows.
#
# Columns removed (by their position in the *original* sheet):
#   B  -> ID / {vendor:id}
#   AA -> Bank Name / {vendor:bank_name}
#   AB -> Bank Account / {vendor:bank_account}
#   AC -> IFSC Code / {vendor:ifsc_code}
#   AD -> Beneficiary Name / {vendor:beneficiary_name}
#   AG -> Brackets Flag / {vendor:brackets_flag}
#
# Delete from right to left so earlier deletions don't shift the column
# letters we still need to reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Brackets Flag (originally column AG)
$ws.Range("AG1").EntireColumn.Delete()

# Bank Name / Bank Account / IFSC Code / Beneficiary Name (originally AA:AD)
$ws.Range("AA1:AD1").EntireColumn.Delete()

# ID (originally column B)
$ws.Range("B1").EntireColumn.Delete()

# Restore the on-screen selection/scroll position to roughly where the
# author left it (now pointing at the SC_Code column area).
$win = $excel.ActiveWindow
$win.ScrollColumn = 23
$ws.Range("AA13").Select()
